# Edit script: add group member names after the title paragraph, and
# merge the "Bryan Nguyen-Le:" / " Question 1" runs into a single run.

$d = $word.ActiveDocument

# --- 1. Insert "Jiacheng Yu", "Hanh Vo", "Bryan Nguyen-Le" as new
#        paragraphs right after the title ("COSC 4353 Software Design")
#        and before the centered "Group 8- Assignment 1" paragraph.
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.Collapse(0)   # wdCollapseEnd
$titleRange.InsertParagraphAfter()
$titleRange.InsertParagraphAfter()
$titleRange.InsertParagraphAfter()

$d.Paragraphs.Item(2).Range.InsertBefore("Jiacheng Yu")
$d.Paragraphs.Item(3).Range.InsertBefore("Hanh Vo")
$d.Paragraphs.Item(4).Range.InsertBefore("Bryan Nguyen-Le")

# --- 2. Merge the two runs that make up the "Bryan Nguyen-Le: Question 1"
#        paragraph into a single run of text (no split between the name
#        and "Question 1").
$d.Content.Find.Execute("Bryan Nguyen-Le: Question 1", $true, $false, $false, $false, $false, $true, 1, $false, "Bryan Nguyen-Le: Question 1", 2)
